# "It now properly uses a calibration data file."
#
# The pH-meter interface workbook's calibration rows were populated with
# placeholder timing values. Update them to the real calibration figures:
#   - "Step" (column D, rows 2-6): 1440 -> 120
#   - "Force delay" (column H, rows 2-6): 5 -> 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D6").Value = 120
$ws.Range("H2:H6").Value = 2

# Leave the selection where the author last clicked.
[void]$ws.Range("H7").Select()
